# OpenTbs 1.8.1 beta - new common keywords for changing cell types in ODS and XLSX.
#
# Renames the "xlsxNum" / "xlsxBool" / "xlsxDate" "ope=" keywords used in the
# demo template to the new common "tbs:num" / "tbs:bool" / "tbs:date" keywords,
# adds a new "Score again" column/example, and adds a named cell ("the_named_cell")
# pointing at a new cell on the "Delete me" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "Examples part 1" ---------------------------------------------
$ws1 = $wb.Worksheets.Item("Examples part 1")

# New "Score again" column next to the existing "Score" column/example.
$ws1.Range("F19").Value = "Score again"
$ws1.Range("E20").Value = "[a.score;ope=tbs:num]"
$ws1.Range("F20").Value = "[a.score;ope=tbs:num]"

# "Merging data with cell" example.
$ws1.Range("C26").Value = "[cell2.score;block=tbs:cell;ope=tbs:num]"

# "Change the type data in a cell" table: xlsxNum/xlsxBool/xlsxDate -> tbs:num/tbs:bool/tbs:date
$ws1.Range("C34").Value = "tbs:num"
$ws1.Range("D34").Value = "[onshow.x_num;ope=tbs:num]"
$ws1.Range("C35").Value = "tbs:bool"
$ws1.Range("D35").Value = "[onshow.x_bt;ope=tbs:bool]"
$ws1.Range("C36").Value = "tbs:date"
$ws1.Range("D36").Value = "[onshow.x_dt;ope=tbs:date]"

# --- Sheet "Delete me" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("Delete me")
$ws4.Range("B6").Value = "And this named cell too."

# --- Workbook defined name pointing at the new cell -----------------------
$wb.Names.Add("the_named_cell", "='Delete me'!`$B`$6")
